# "minor adjustment su rope"
#
# Updates the "Sequenza effettiva" (actual lesson sequence) notes for the
# last few lessons, turns the E15:E20 date formulas into one shared
# formula group, adds a new note in F21, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of the "effective lesson" column (F) for the last rows
# of the syllabus to reflect what was actually covered.
$ws.Range("F16").Value = "chiuso MCMC, prob programming (fino a cambiare prior beta(1,1) con uniform)"
$ws.Range("F18").Value = "prob programming (finire nbook 1 con rope);  ultimi 15 min presenta assignment "
$ws.Range("F19").Value = "normal-normal (ultimi 15 min presentare progetto)"
$ws.Range("F20").Value = "hyp test"

# Row 21 previously had no note in column F; add a placeholder note.
$ws.Range("F21").Value = " "

# Re-enter the "data" formulas for E15:E20 as one range assignment so they
# become a single shared formula group (same relative formula E[n-2]+7).
$ws.Range("E15:E20").Formula = "=E13+7"

# Move the active selection to F18, matching the latest edit location.
$ws.Range("F18").Select()
